# Reverted program to produce separate excel files
#
# This workbook is a flattened ("pasted as values") export of several
# small pivot-table summaries stacked vertically on one sheet. The edit:
#   1. Renames the sheet "PivotTable" -> "Pivot Table".
#   2. Relabels a status value "ffxqa" -> "open" wherever it is used as a
#      row label.
#   3. Relabels the status categories "Customer Testing" -> "Customer
#      Generated Ticket" and "FFX Testing" -> "Customer Testing".
#   4. Updates ~75 recomputed counts/subtotals across the six stacked
#      pivot summaries to match the relabelled/regrouped data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Sheet name -----------------------------------------------------
$ws.Name = "Pivot Table"

# --- 2. Status label "ffxqa" -> "open" ---------------------------------
$ws.Range("A10").Value = "open"
$ws.Range("A19").Value = "open"

# --- 3. Status category relabels ---------------------------------------
# "Customer Testing" -> "Customer Generated Ticket"
$ws.Range("C2").Value = "Customer Generated Ticket"
$ws.Range("C18").Value = "Customer Generated Ticket"
$ws.Range("A57").Value = "Customer Generated Ticket"

# "FFX Testing" -> "Customer Testing"
$ws.Range("D2").Value = "Customer Testing"
$ws.Range("D18").Value = "Customer Testing"
$ws.Range("A58").Value = "Customer Testing"

# --- 4. Recomputed numeric values ---------------------------------------
# Table 1 (rows 2-6): Count of CAC/MOF Requestor
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = 6

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("F5").Value = 6

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 3
$ws.Range("F6").Value = 13

# Table 2 (rows 8-14): Count of Status
$ws.Range("B10").Value = 0
$ws.Range("D10").Value = 1

$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 3

$ws.Range("B14").Value = 6
$ws.Range("C14").Value = 4
$ws.Range("E14").Value = 13

# Table 3 (rows 17-23): Count of Status
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0

$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 3

$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 3
$ws.Range("F23").Value = 13

# Table 4 (rows 26-31): Count of CAC/MOF/FFX Owner
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 3

$ws.Range("B29").Value = 2
$ws.Range("D29").Value = 1

$ws.Range("C30").Value = 2
$ws.Range("E30").Value = 5

$ws.Range("B31").Value = 6
$ws.Range("C31").Value = 4
$ws.Range("E31").Value = 13

# Table 5 (rows 33-38): Count of CAC/MOF/FFX Owner (by Severity)
$ws.Range("B35").Value = 2
$ws.Range("D35").Value = 3
$ws.Range("E35").Value = 6

$ws.Range("B36").Value = 2
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = 4

$ws.Range("B38").Value = 6
$ws.Range("D38").Value = 6
$ws.Range("E38").Value = 13

# Table 6 (rows 40-44): Count of Ministry/FFX Owner
$ws.Range("B42").Value = 2
$ws.Range("D42").Value = 1

$ws.Range("B43").Value = 4
$ws.Range("C43").Value = 2
$ws.Range("D43").Value = 2
$ws.Range("E43").Value = 8

$ws.Range("B44").Value = 6
$ws.Range("C44").Value = 4
$ws.Range("E44").Value = 13

# Table 7 (rows 47-52): Count of CAC/MOF/FFX Owner
$ws.Range("B49").Value = 2
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 1
$ws.Range("E49").Value = 3

$ws.Range("B50").Value = 2
$ws.Range("D50").Value = 1

$ws.Range("C51").Value = 2
$ws.Range("E51").Value = 5

$ws.Range("B52").Value = 6
$ws.Range("C52").Value = 4
$ws.Range("E52").Value = 13

# Table 8 (rows 54-60): Count of Status
$ws.Range("B57").Value = 0
$ws.Range("C57").Value = 0
$ws.Range("E57").Value = 1

$ws.Range("B58").Value = 2
$ws.Range("C58").Value = 1
$ws.Range("E58").Value = 3

$ws.Range("B60").Value = 6
$ws.Range("C60").Value = 4
$ws.Range("E60").Value = 13
